$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the existing weekly rows down.
$ws.Rows.Item(2).Insert()

# New row 2 holds this week's freshly published data; reuse the same
# layout/values as the rest of the table, just with the new figures.
$ws.Cells.Item(2, 1).Value = 8
$ws.Cells.Item(2, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(2, 3).Value = "Coquimbo"
$ws.Cells.Item(2, 4).Value = 44881
$ws.Cells.Item(2, 4).Style = $ws.Cells.Item(3, 4).Style
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 100112039
$ws.Cells.Item(2, 7).Value = "Ciboulette"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 500
$ws.Cells.Item(2, 11).Value = 1900
$ws.Cells.Item(2, 12).Value = 2000
$ws.Cells.Item(2, 13).Value = 1950
$ws.Cells.Item(2, 14).Value = "`$/docena de atados"
$ws.Cells.Item(2, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(2, 16).Value = 650
$ws.Cells.Item(2, 17).Value = 3
$ws.Cells.Item(2, 18).Value = "Hortaliza"

# Clear the bold/border formatting that Insert() copied down from the
# header row onto the non-date cells of the new row.
$headerStyleCols = @(1,2,3,5,6,7,8,9,10,11,12,13,14,15,16,17,18)
foreach ($col in $headerStyleCols) {
    $ws.Cells.Item(2, $col).Style = $ws.Cells.Item(3, $col).Style
}
